$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "password123"
$ws.Range("A3").Value = "admin123"
$ws.Range("B3").Value = "admin321"
$ws.Range("A4").Value = "admin123"
$ws.Range("B4").Value = "admin321"
$ws.Range("A5").Value = "admin123"
$ws.Range("B5").Value = "admin321"
$ws.Range("A6").Value = "admin123"
$ws.Range("B6").Value = "admin321"
